# Refresh the "cryptos" price/volume table with the latest scraped values.
# A couple of Price cells (D18, D19, D39, D41) are numeric-looking strings
# that end in a trailing zero (e.g. "69.00", "0.0190"); Excel would silently
# normalise those to plain numbers (69, 0.019) if assigned as-is, so they are
# written with a leading apostrophe to force them to stay literal text,
# exactly like typing '69.00 into the cell by hand.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.479.02'
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("D3").Value = '1.795.10'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '223.02'
$ws.Range("E5").Value = '  -1.80%  '
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '32.25'
$ws.Range("E8").Value = '  +1.79%  '
$ws.Range("D9").Value = '0.288'
$ws.Range("E9").Value = '  +2.32%  '
$ws.Range("D10").Value = '0.0705'
$ws.Range("E10").Value = '  +6.83%  '
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").Value = '2.053.71'
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("D13").Value = '10.96'
$ws.Range("E13").Value = '  -4.49%  '
$ws.Range("D14").Value = '1.784.19'
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").Value = '0.638'
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").Value = '34.490.10'
$ws.Range("E16").Value = '  +1.10%  '
$ws.Range("D17").Value = '4.27'
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("D18").Value = '''69.00'
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("D19").Value = '''250.60'
$ws.Range("E19").Value = '  -1.50%  '
$ws.Range("D20").Value = '0.0₃0796'
$ws.Range("E20").Value = '  +7.09%  '
$ws.Range("D21").Value = '11.02'
$ws.Range("E21").Value = '  +4.94%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("E23").Value = '  -0.86%  '
$ws.Range("D24").Value = '2.16'
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("D25").Value = '161.33'
$ws.Range("E25").Value = '  +2.62%  '
$ws.Range("D26").Value = '16.34'
$ws.Range("E26").Value = '  -1.58%  '
$ws.Range("D27").Value = '7.12'
$ws.Range("E27").Value = '  +0.99%  '
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("D30").Value = '552.56'
$ws.Range("E30").Value = '  +963.95%  '
$ws.Range("D31").Value = '0.0523'
$ws.Range("E31").Value = '  +1.09%  '
$ws.Range("D32").Value = '3.77'
$ws.Range("E32").Value = '  -1.51%  '
$ws.Range("E33").Value = '  -1.10%  '
$ws.Range("D34").Value = '3.58'
$ws.Range("E34").Value = '  -0.66%  '
$ws.Range("D35").Value = '1.88'
$ws.Range("E35").Value = '  +1.73%  '
$ws.Range("D36").Value = '1.422.59'
$ws.Range("E36").Value = '  -1.91%  '
$ws.Range("E37").Value = '  -0.52%  '
$ws.Range("D38").Value = '0.636'
$ws.Range("E38").Value = '  +0.12%  '
$ws.Range("D39").Value = '''0.0190'
$ws.Range("E39").Value = '  +1.06%  '
$ws.Range("D40").Value = '82.97'
$ws.Range("E40").Value = '  -0.81%  '
$ws.Range("D41").Value = '''0.950'
$ws.Range("E41").Value = '  +5.34%  '
$ws.Range("E42").Value = '  -3.55%  '
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("D44").Value = '2.12'
$ws.Range("E44").Value = '  +1.53%  '
$ws.Range("D45").Value = '6.02'
$ws.Range("E45").Value = '  +2.32%  '
$ws.Range("E46").Value = '  -0.95%  '
$ws.Range("E47").Value = '  -2.50%  '
$ws.Range("D48").Value = '1.946.11'
$ws.Range("E48").Value = '  +0.26%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '12.21'
$ws.Range("E49").Value = '  +0.62%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '105.58'
$ws.Range("E50").Value = '  +7.35%  '
$ws.Range("E51").Value = '  -0.05%  '
